$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.347467
$ws.Range("H2").Value = 1.042401
$ws.Range("I2").Value = 0.864291305025554
$ws.Range("J2").Value = 0.864291305025554
$ws.Range("M2").Value = 52.47402833333333
$ws.Range("N2").Value = 157.422085
$ws.Range("O2").Value = 0.699720168977827
$ws.Range("P2").Value = 0.6997201689778269
$ws.Range("Q2").Value = 18.23299320289833
$ws.Range("R2").Value = 164.096938826085
$ws.Range("S2").Value = 0.6047620579985472
$ws.Range("T2").Value = 0.6047620579985471

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.347467
$ws.Range("H3").Value = 1.042401
$ws.Range("I3").Value = 0.864291305025554
$ws.Range("J3").Value = 0.864291305025554
$ws.Range("O3").Value = 0.1650338345468634
$ws.Range("P3").Value = 0.1650338345468634
$ws.Range("Q3").Value = 4.300377375054
$ws.Range("R3").Value = 38.703396375486
$ws.Range("S3").Value = 0.14263730823388
$ws.Range("T3").Value = 0.14263730823388

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.347467
$ws.Range("H4").Value = 1.042401
$ws.Range("I4").Value = 0.864291305025554
$ws.Range("J4").Value = 0.864291305025554
$ws.Range("M4").Value = 6.377905999999999
$ws.Range("N4").Value = 19.133718
$ws.Range("O4").Value = 0.08504682422503862
$ws.Range("P4").Value = 0.08504682422503862
$ws.Range("Q4").Value = 2.216111864101999
$ws.Range("R4").Value = 19.945006776918
$ws.Range("S4").Value = 0.07350523069773753
$ws.Range("T4").Value = 0.07350523069773753

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.347467
$ws.Range("H5").Value = 1.042401
$ws.Range("I5").Value = 0.864291305025554
$ws.Range("J5").Value = 0.864291305025554
$ws.Range("M5").Value = 3.764580333333333
$ws.Range("N5").Value = 11.293741
$ws.Range("O5").Value = 0.05019917225027106
$ws.Range("P5").Value = 0.05019917225027107
$ws.Range("Q5").Value = 1.308067434682333
$ws.Range("R5").Value = 11.772606912141
$ws.Range("S5").Value = 0.04338670809538935
$ws.Range("T5").Value = 0.04338670809538935

# Row 6
$ws.Range("I6").Value = 0.135708694974446
$ws.Range("J6").Value = 0.135708694974446
$ws.Range("M6").Value = 52.47402833333333
$ws.Range("N6").Value = 157.422085
$ws.Range("O6").Value = 0.699720168977827
$ws.Range("P6").Value = 0.6997201689778269
$ws.Range("Q6").Value = 2.862895529152777
$ws.Range("R6").Value = 25.76605976237499
$ws.Range("S6").Value = 0.09495811097927977
$ws.Range("T6").Value = 0.09495811097927975

# Row 7
$ws.Range("I7").Value = 0.135708694974446
$ws.Range("J7").Value = 0.135708694974446
$ws.Range("O7").Value = 0.1650338345468634
$ws.Range("P7").Value = 0.1650338345468634
$ws.Range("S7").Value = 0.02239652631298349
$ws.Range("T7").Value = 0.02239652631298349

# Row 8
$ws.Range("I8").Value = 0.135708694974446
$ws.Range("J8").Value = 0.135708694974446
$ws.Range("M8").Value = 6.377905999999999
$ws.Range("N8").Value = 19.133718
$ws.Range("O8").Value = 0.08504682422503862
$ws.Range("P8").Value = 0.08504682422503862
$ws.Range("Q8").Value = 0.3479679215166666
$ws.Range("R8").Value = 3.13171129365
$ws.Range("S8").Value = 0.01154159352730109
$ws.Range("T8").Value = 0.01154159352730109

# Row 9
$ws.Range("I9").Value = 0.135708694974446
$ws.Range("J9").Value = 0.135708694974446
$ws.Range("M9").Value = 3.764580333333333
$ws.Range("N9").Value = 11.293741
$ws.Range("O9").Value = 0.05019917225027106
$ws.Range("P9").Value = 0.05019917225027107
$ws.Range("Q9").Value = 0.2053892286861111
$ws.Range("R9").Value = 1.848503058175
$ws.Range("S9").Value = 0.006812464154881711
$ws.Range("T9").Value = 0.006812464154881712
